$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 -> new values (becomes the former row 5 content, but with own B value)
$ws.Range("A3").Value = 112105682
$ws.Range("B3").Value = 89104
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 5747
$ws.Range("F3").Value = "Läderdoftande fingersvamp"
$ws.Range("G3").Value = "Ramaria safraniolens"
$ws.Range("H3").Value = "Christian"
$ws.Range("P3").Value = "Svensbergsbäcken (Svensbergsbäcken), Jmt"
$ws.Range("Q3").Value = 446627
$ws.Range("R3").Value = 7032919

# Row 4 -> only taxonsorteringsordning changes
$ws.Range("B4").Value = 84943

# Row 5 -> new values (becomes the former row 3 content, but with own B value)
$ws.Range("A5").Value = 112105307
$ws.Range("B5").Value = 89114
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5754
$ws.Range("F5").Value = "Gultoppig fingersvamp"
$ws.Range("G5").Value = "Ramaria testaceoflava"
$ws.Range("H5").Value = "(Bres.) Corner"
$ws.Range("P5").Value = "Landverktjärnen (Landverktjärnen), Jmt"
$ws.Range("Q5").Value = 446544
$ws.Range("R5").Value = 7032738

# Row 6 -> only taxonsorteringsordning changes
$ws.Range("B6").Value = 90835
